$wb = $excel.ActiveWorkbook

$wsVar = $wb.Worksheets.Item("Variables")
$wsCat = $wb.Worksheets.Item("Categories")

# --- Sheet "Variables": add new row 109 for the DIETARY_ASSESS_INSTR variable ---
# Write the index/name/type first so the "DIETARY_ASSESS_INSTR" shared string is
# registered before the Categories sheet entries (matches original authoring order).
$wsVar.Range("A109").Value = 108
$wsVar.Range("B109").Value = "DIETARY_ASSESS_INSTR"
$wsVar.Range("D109").Value = "integer"

# --- Sheet "Categories": add the 7 category rows for DIETARY_ASSESS_INSTR ---
$wsCat.Range("A89").Value = "DIETARY_ASSESS_INSTR"
$wsCat.Range("B89").Value = "FPQ (Food propensity questionnair without portion sizes)"
$wsCat.Range("C89").Value = 0

$wsCat.Range("A90").Value = "DIETARY_ASSESS_INSTR"
$wsCat.Range("B90").Value = "FFQ (Food frequency questionnaire"
$wsCat.Range("C90").Value = 1

$wsCat.Range("A91").Value = "DIETARY_ASSESS_INSTR"
$wsCat.Range("B91").Value = "24HDR (24-h dietary recall"
$wsCat.Range("C91").Value = 2

$wsCat.Range("A92").Value = "DIETARY_ASSESS_INSTR"
$wsCat.Range("B92").Value = "3_d_FR_w (3-day weighing food record)"
$wsCat.Range("C92").Value = 3

$wsCat.Range("A93").Value = "DIETARY_ASSESS_INSTR"
$wsCat.Range("B93").Value = "7_d_FR (7-day  food record; described portion sizes)"
$wsCat.Range("C93").Value = 4

$wsCat.Range("A94").Value = "DIETARY_ASSESS_INSTR"
$wsCat.Range("B94").Value = "7_d_FR_w (7-day weighing food record)"
$wsCat.Range("C94").Value = 5

$wsCat.Range("A95").Value = "DIETARY_ASSESS_INSTR"
$wsCat.Range("B95").Value = "24HFL_FFQ (24-h short food list combined with FFQ"
$wsCat.Range("C95").Value = 6

# --- Back to "Variables": finish row 109 with its label (added last so this
# string lands after the Categories strings in the shared string table) ---
$wsVar.Range("C109").Value = "Dietary Assessment Instrument"

# --- Update selection / scroll position to mirror final cursor location ---
$wsCat.Activate()
$wsCat.Range("A89:C95").Select()

$wsVar.Activate()
$wsVar.Range("A110").Select()
